$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 92; this pushes the former
# rows 92-118 down to 94-120 (weekly price list data shifts down
# to make room for the two new "Early Burlat" entries at the top
# of this week's batch).
$ws.Range("A92:A93").EntireRow.Insert()

# New row 92: Early Burlat / Primera
$ws.Range("A92").Value = 10
$ws.Range("B92").Value = "Vega Modelo de Temuco"
$ws.Range("C92").Value = "La Araucanía"
$ws.Range("D92").Value = 44511
$ws.Range("E92").Value = 9
$ws.Range("F92").Value = "Fruta"
$ws.Range("G92").Value = 100103
$ws.Range("H92").Value = "Frutos de hueso (carozo)"
$ws.Range("I92").Value = 100103001
$ws.Range("J92").Value = "Cereza"
$ws.Range("K92").Value = "Early Burlat"
$ws.Range("L92").Value = "Primera"
$ws.Range("M92").Value = 95
$ws.Range("N92").Value = 21000
$ws.Range("O92").Value = 21000
$ws.Range("P92").Value = 21000
$ws.Range("Q92").Value = "$/bandeja 6 kilos"
$ws.Range("R92").Value = "Región Metropolitana"
$ws.Range("S92").Value = 3500
$ws.Range("T92").Value = 6

# New row 93: Early Burlat / Primera
$ws.Range("A93").Value = 10
$ws.Range("B93").Value = "Vega Modelo de Temuco"
$ws.Range("C93").Value = "La Araucanía"
$ws.Range("D93").Value = 44511
$ws.Range("E93").Value = 9
$ws.Range("F93").Value = "Fruta"
$ws.Range("G93").Value = 100103
$ws.Range("H93").Value = "Frutos de hueso (carozo)"
$ws.Range("I93").Value = 100103001
$ws.Range("J93").Value = "Cereza"
$ws.Range("K93").Value = "Early Burlat"
$ws.Range("L93").Value = "Primera"
$ws.Range("M93").Value = 110
$ws.Range("N93").Value = 3500
$ws.Range("O93").Value = 3500
$ws.Range("P93").Value = 3500
$ws.Range("Q93").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R93").Value = "Región del Maule"
$ws.Range("S93").Value = 3500
$ws.Range("T93").Value = 1
